$d = $word.ActiveDocument

# 1) Update the two outer paragraphs' text.
$d.Content.Find.Execute("Outer para 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "outer, before sect break", 2)
$d.Content.Find.Execute("Outer para 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "outer, after sect break", 2)

# 2) Insert a (continuous) section break right after paragraph 1, then
#    delete the now-redundant extra paragraph mark so the new sectPr
#    ends up living inside paragraph 1's own pPr (matching how the
#    fixture wants a sectPr-bearing paragraph without adding a blank
#    paragraph).
$p1 = $d.Paragraphs(1)
$breakPos = $p1.Range.End - 1
$r = $d.Range($breakPos, $breakPos)
$r.InsertBreak(3)
$d.Range($breakPos, $breakPos + 1).Delete()

# 3) Fix up page geometry on every section: Letter -> A4, 720 -> 708 twips
#    for header/footer distance and column spacing.
for ($i = 1; $i -le $d.Sections.Count; $i++) {
    $s = $d.Sections($i)
    $s.PageSetup.PageWidth = 595.3
    $s.PageSetup.PageHeight = 841.9
    $s.PageSetup.HeaderDistance = 35.4
    $s.PageSetup.FooterDistance = 35.4
    $s.PageSetup.TextColumns.Spacing = 35.4
}

# 4) First section restarts page numbering at 1.
$d.Sections(1).Headers(1).PageNumbers.StartingNumber = 1
